# Update functions and Data Model (#50)
# - Fix "Auteur·ice" -> "Auteur·rice" in D2
# - Add three new properties (rows 38-40): hasAuthorshipResource, hasCopyrightResource,
#   hasLicenseResource, mirroring the existing comment_* columns on the en/de/fr/it labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the French label for "Authorship" (typo correction: Auteur·ice -> Auteur·rice)
$ws.Range("D2").Value = "Auteur·rice"

# 2. New row 38: hasAuthorshipResource
$ws.Range("A38").Value = "hasAuthorshipResource"
$ws.Range("B38").Value = "Author of the resource"
$ws.Range("C38").Value = "Autor der Resource"
$ws.Range("D38").Value = "Auteur·rice de la ressource"
$ws.Range("E38").Value = "Autore della risorsa"
$ws.Range("G38").Value = "Author of the resource"
$ws.Range("H38").Value = "Autor der Resource"
$ws.Range("I38").Value = "Auteur·rice de la ressource"
$ws.Range("J38").Value = "Autore della risorsa"
$ws.Range("L38").Value = "hasValue, foaf:person, schema:author, crm:E21_Person"
$ws.Range("M38").Value = "TextValue"
$ws.Range("N38").Value = "SimpleText"

# 3. New row 39: hasCopyrightResource
$ws.Range("A39").Value = "hasCopyrightResource"
$ws.Range("B39").Value = "Copyright of the resource"
$ws.Range("C39").Value = "Urheberrecht der Resource"
$ws.Range("D39").Value = "Droits d'auteur de la ressource"
$ws.Range("E39").Value = "Copyright della risorsa"
$ws.Range("G39").Value = "Copyright of the resource"
$ws.Range("H39").Value = "Urheberrecht der Resource"
$ws.Range("I39").Value = "Droits d'auteur de la ressource"
$ws.Range("J39").Value = "Copyright della risorsa"
$ws.Range("L39").Value = "hasValue, schema:copyrightHolder, crm:P105_right_held_by"
$ws.Range("M39").Value = "TextValue"
$ws.Range("N39").Value = "SimpleText"
$ws.Range("A39").Font.Color = 3355443

# 4. New row 40: hasLicenseResource
$ws.Range("A40").Value = "hasLicenseResource"
$ws.Range("B40").Value = "License of the resource"
$ws.Range("C40").Value = "Lizenz der Resource"
$ws.Range("D40").Value = "Licence de la ressource"
$ws.Range("E40").Value = "Licenza della risorsa"
$ws.Range("G40").Value = "License of the resource"
$ws.Range("H40").Value = "Lizenz der Resource"
$ws.Range("I40").Value = "Licence de la ressource"
$ws.Range("J40").Value = "Licenza della risorsa"
$ws.Range("L40").Value = "hasValue, schema:license"
$ws.Range("M40").Value = "ListValue"
$ws.Range("N40").Value = "List"
$ws.Range("O40").Value = "hlist: License"
$ws.Range("A40").Font.Color = 3355443

# 5. Match the author's final selection / cursor position
$ws.Range("L40").Select()
